# "hr portion added" - append a new employee record (row 2) to the
# "Personal" sheet of the employee-import template, including an
# Email hyperlink, matching the author's commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personal")

$ws.Range("A2").Value = "NGO0001002"
$ws.Range("B2").Value = "Nathan"
$ws.Range("C2").Value = "McCullam"
$ws.Range("D2").Value = "Male"
$ws.Range("E2").Value = "Christianity"
$ws.Range("F2").Value = "B+"
$ws.Range("G2").Value = "Separated"
$ws.Range("H2").Value = "10/10/1996"
$ws.Range("I2").Value = "Bangladeshi"
$ws.Range("J2").Value = "Nathan Gilbert"
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 1720906313
$ws.Range("N2").Value = "noyon9718nil@gmail.com"
$ws.Range("O2").Value = "NID"
$ws.Range("P2").Value = 90284828439997
$ws.Range("Q2").Value = "Chadpur"
$ws.Range("R2").Value = "Dhaka"

# Turn the email into a real mailto hyperlink (adds font 1 / "Hyperlink"
# cell style automatically, like Excel does).
$ws.Hyperlinks.Add($ws.Range("N2"), "mailto:noyon9718nil@gmail.com", "", "", "noyon9718nil@gmail.com") | Out-Null

# Widen the Mobile / Email columns that now hold longer values.
$ws.Columns.Item(13).ColumnWidth = 10.166666666666666
$ws.Columns.Item(14).ColumnWidth = 23
$ws.Columns.Item(16).ColumnWidth = 11.166666666666666

# Leave the selection / scroll position where the author left it.
$ws.Range("Q3").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 7
